$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "benthic-methane cycling"
$ws.Range("A4").Value = "benthic-non-methane cycling"
$ws.Range("B4").Value = -30
$ws.Range("C4").Value = -22

$ws.Range("A5").Select()
